$wb = $excel.ActiveWorkbook

# 1) Rename the "Include" sheet tab.
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# 2) Update the Metadata sheet.
$ws = $wb.Worksheets.Item(1)

# Update the Date value (row 8, column B).
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row above row 11 ("Description") for the new "Jurisdiction" property,
# copying the formatting from the row above so the new cells keep the same style.
$ws.Range("A11").EntireRow.Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
